$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-14 (columns D, J, K, L, M, N, O, P, Q).
# The data rows were reordered/shuffled; other columns (A,B,C,E,F,G,H,I,R)
# stay the same for every row, so only these columns need updating.
$rows = @(
    @{ Row=2;  D=44209; J=150; K=3500; L=4000; M=3767; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1884; Q=2 },
    @{ Row=3;  D=44701; J=120; K=7000; L=7500; M=7250; N='$/paquete 36 unidades'; O='Región Metropolitana';   P=201;  Q=36 },
    @{ Row=4;  D=44215; J=140; K=3500; L=4000; M=3768; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1884; Q=2 },
    @{ Row=5;  D=44161; J=50;  K=2800; L=3000; M=2900; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1450; Q=2 },
    @{ Row=6;  D=44210; J=105; K=3500; L=4000; M=3714; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1857; Q=2 },
    @{ Row=7;  D=44225; J=80;  K=3400; L=3700; M=3550; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1775; Q=2 },
    @{ Row=8;  D=44662; J=200; K=8000; L=8500; M=8250; N='$/paquete 36 unidades'; O='Región Metropolitana';   P=229;  Q=36 },
    @{ Row=9;  D=44223; J=80;  K=3500; L=3800; M=3688; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1844; Q=2 },
    @{ Row=10; D=44208; J=85;  K=3700; L=4000; M=3824; N='$/paquete 2 kilos';    O='Provincia de Diguillín'; P=1912; Q=2 },
    @{ Row=11; D=44664; J=200; K=8000; L=8500; M=8250; N='$/paquete 36 unidades'; O='Región Metropolitana';   P=229;  Q=36 },
    @{ Row=12; D=44160; J=43;  K=3500; L=4000; M=3709; N='$/paquete 36 unidades'; O='Región Metropolitana';   P=103;  Q=36 },
    @{ Row=13; D=44166; J=70;  K=3500; L=4000; M=3679; N='$/paquete 36 unidades'; O='Región Metropolitana';   P=102;  Q=36 },
    @{ Row=14; D=44704; J=100; K=6000; L=6500; M=6250; N='$/paquete 36 unidades'; O='Región Metropolitana';   P=174;  Q=36 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
}
